# Applies the pedalboard-led-ring-bom.xlsx BoM refresh:
#   - capacitor footprint C_0201_0603Metric_Pad0.64x0.40mm_HandSolder -> C_0402_1005Metric
#   - connector value "middle" -> "02x02"
#   - resistor footprint R_0201_0603Metric_Pad0.64x0.40mm_HandSolder -> R_0402_1005Metric
#   - "Created:" timestamp 2023-10-20 09:27:55 -> 2023-10-21 09:13:57
#   - BoM sheet Footprint column (F) narrower: 48.7109375 -> 40.7109375 chars
#   - Costs sheet Footprint column (C) narrower: 44.7109375 -> 36.7109375 chars

$wb = $excel.ActiveWorkbook

$bom = $wb.Worksheets.Item("BoM")
$bom.Range("F9").Value = "C_0402_1005Metric"
$bom.Range("E11").Value = "02x02"
$bom.Range("F12").Value = "R_0402_1005Metric"

$costs = $wb.Worksheets.Item("Costs")
$costs.Range("C10").Value = "C_0402_1005Metric"
$costs.Range("B12").Value = "02x02"
$costs.Range("C13").Value = "R_0402_1005Metric"
$costs.Range("B16").Value = "2023-10-21 09:13:57"

$costsDnf = $wb.Worksheets.Item("Costs (DNF)")
$costsDnf.Range("B13").Value = "2023-10-21 09:13:57"

# Column width shrinks (character widths correspond to the stored
# 40.7109375 / 36.7109375 sheet XML widths).
$bom.Columns.Item(6).ColumnWidth = 39.83
$costs.Columns.Item(3).ColumnWidth = 35.83
